$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.369.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.90%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.775.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.63%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.27%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.17%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'306.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.95%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4227"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.3599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.58%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.52%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8354"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.07%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'20.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.46%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.783.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.34%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'6.452"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.63%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.238"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.03%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.06866"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.51%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.23%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'78.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.35%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008634"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.74%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.00%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.19%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'26.370.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.97%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.082"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.05%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.10%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.000.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.47%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'152.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.97%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.809"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -7.70%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'17.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.68%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'5.063"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.49%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'114.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.71%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.829"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +9.47%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.08851"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.37%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.7263"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.11%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.118"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.39%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.318"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.39%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.000"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.25%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.732"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -7.02%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.091"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.01%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.05128"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.19%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01883"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.71%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.4917"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.58%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1608"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.81%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -3.62%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'6.326"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.44%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'7.963"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.47%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'104.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.06%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.15%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'10.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.12%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.633"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.76%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.06170"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.75%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.4440"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.38%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.729"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +4.23%  "
$ws.Range("E51").Style = "Normal"
